$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transmittals_New")

# The M column formulas (rows 3-11) all referenced N2 instead of the row's own
# N cell. Repoint each formula to its own row so the "message" text correctly
# reflects Delegate/Reply-All overrides in N8:N11 (and keeps parity for the
# unaffected rows 3-7, which still resolve to "Message for New transmittal").
for ($row = 3; $row -le 11; $row++) {
    $ws.Range("M$row").Formula = "=CONCATENATE(ROW()-1,`" of `",COUNTA(A2:A100),`" `",N$row)"
}
